$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("testData")

# Delete row 2 (mngr276899/qapydAq), shifting rows 3-5 up to 2-4
$ws.Rows.Item(2).Delete()

# Update selection to match target state
$ws.Range("A2:B2").Select()
